# Replace the digit strings in column A (English) and column B (French)
# for the "numbers" section of the vocabulary list (rows 374-406) with
# their word equivalents. Column C (Fulfulde) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$english = @(
    "one", "two ", "three", "four", "five", "six", "seven", "eight", "nine", "ten",
    "elleven", "twelve", "thirteen", "fourteen", "fifteen", "sixteen", "seventeen", "eighteen", "nineteen", "twenty",
    "twenty-one", "twenty-two", "twenty-three", "thirty", "forty", "fifty", "sixty", "seventy", "eighty", "ninety",
    "one-hundred", "two-hundred", "three-hundred"
)

$french = @(
    "un", "deux", "trois", "quatre", "cinq", "six", "sept", "huit", "neuf", "dix",
    "onze", "douze", "treize", "quatorze", "quinze", "seize", "dix-sept", "dix-huit", "dix-neuf", "vingt",
    "vingt et un", "vingt-deux", "vingt-trois", "trente", "quarante", "cinquante", "soixante", "soixante-dix", "quatre-vingts", "quatre-vingt-dix",
    "cent", "deux cents", "trois cents"
)

$startRow = 374

for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $english[$i]
    $ws.Cells.Item($row, 2).Value = $french[$i]
}

$ws.Range("E405").Select()
